$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (E1, "panel") onto the
# new header cell F1, then set its text — mirrors the s="1" style shared
# by all header cells in the sheet.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:41:37.745622",
    "2021-10-05 13:41:37.745636",
    "2021-10-05 13:41:37.745640",
    "2021-10-05 13:41:37.745643",
    "2021-10-05 13:41:37.745646",
    "2021-10-05 13:41:37.745650",
    "2021-10-05 13:41:37.745653",
    "2021-10-05 13:41:37.745656",
    "2021-10-05 13:41:37.745659",
    "2021-10-05 13:41:37.745662",
    "2021-10-05 13:41:37.745665",
    "2021-10-05 13:41:37.745668",
    "2021-10-05 13:41:37.745671",
    "2021-10-05 13:41:37.745674",
    "2021-10-05 13:41:37.745677",
    "2021-10-05 13:41:37.745680",
    "2021-10-05 13:41:37.745683",
    "2021-10-05 13:41:37.745686",
    "2021-10-05 13:41:37.745689",
    "2021-10-05 13:41:37.745692",
    "2021-10-05 13:41:37.745695",
    "2021-10-05 13:41:37.745698",
    "2021-10-05 13:41:37.745701",
    "2021-10-05 13:41:37.745704",
    "2021-10-05 13:41:37.745708",
    "2021-10-05 13:41:37.745711",
    "2021-10-05 13:41:37.745714",
    "2021-10-05 13:41:37.745717",
    "2021-10-05 13:41:37.745720",
    "2021-10-05 13:41:37.745723"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
